$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated "Price" values are plain decimal numbers (e.g. "564.98") that Excel
# would otherwise auto-convert to a numeric cell. The source data keeps these as
# plain text, so force a text format before assigning, then clear the temporary
# formatting again so the cells end up with no explicit style, just like before.
$textCells = @("D5", "D6", "D12", "D17", "D19", "D21", "D22", "D23", "D25", "D26", "D28", "D32", "D33", "D34", "D36", "D37", "D39", "D42", "D44", "D45", "D47")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.162.35"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").Value = "2.484.34"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "564.98"
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("D6").Value = "163.36"
$ws.Range("E6").Value = "  -4.38%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("D9").Value = "2.479.60"
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "0.352"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "2.937.33"
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("D15").Value = "69.036.44"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("E16").Value = "  -3.10%  "
$ws.Range("D17").Value = "24.17"
$ws.Range("E17").Value = "  -5.08%  "
$ws.Range("D18").Value = "2.474.16"
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("D19").Value = "11.13"
$ws.Range("E19").Value = "  -4.07%  "
$ws.Range("E20").Value = "  -7.98%  "
$ws.Range("D21").Value = "343.53"
$ws.Range("E21").Value = "  -3.55%  "
$ws.Range("D22").Value = "3.85"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("D23").Value = "1.91"
$ws.Range("E23").Value = "  -7.33%  "
$ws.Range("D25").Value = "69.16"
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -5.09%  "
$ws.Range("D27").Value = "2.608.65"
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("D28").Value = "8.64"
$ws.Range("E28").Value = "  -5.30%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "0.0₃0868"
$ws.Range("E30").Value = "  -5.80%  "
$ws.Range("E31").Value = "  -4.20%  "
$ws.Range("D32").Value = "440.16"
$ws.Range("E32").Value = "  -6.49%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.18"
$ws.Range("E33").Value = "  -7.99%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  -4.29%  "
$ws.Range("D36").Value = "155.22"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").Value = "0.113"
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").Value = "18.07"
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").Value = "4.56"
$ws.Range("E42").Value = "  -6.38%  "
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("D44").Value = "37.86"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "2.16"
$ws.Range("E45").Value = "  -9.07%  "
$ws.Range("E46").Value = "  -8.94%  "
$ws.Range("D47").Value = "138.08"
$ws.Range("E47").Value = "  -4.99%  "
$ws.Range("E48").Value = "  -3.75%  "
$ws.Range("E49").Value = "  -5.35%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("E51").Value = "  -2.49%  "

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).ClearFormats()
}
